$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E1 header (text unchanged, but triggers the shared-string table rebuild)
$ws.Range("E1").Value = "time_per_use_min"

$letters = @("A", "B", "C", "D", "E", "F", "G", "H", "J", "K", "L", "M", "N", "O")
$names   = @("Nevera", "Iluminación", "Ducha", "TV", "Codificador TV", "Plancha", "Lavadora", "PC", "Laptop", "Cargadores de celular", "Microondas", "Licuadora", "Router", "Equipo Sonido")
$power   = @(120, 30, 3500, 70, 80, 1100, 400, 140, 60, $null, 1250, 400, 20, 50)
$weekly  = @(42, 49, 28, 3.5, 6, 20, 1, 1, 1, 1, 1, 1, 1, 1)
$timeuse = @(30, 30, 20, 10, 3, 5, 1, 1, 1, 1, 1, 1, 1, 1)

# Fill column A (ids / letters) for all rows first
for ($i = 0; $i -lt $letters.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $letters[$i]
}

# Then fill column B (device names) for all rows -- last row's name first,
# matching the author's original entry order, then the remaining rows in order
$ws.Cells.Item(15, 2).Value = $names[13]
for ($i = 0; $i -lt $names.Count - 1; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $names[$i]
}

# Then columns C, D, E
for ($i = 0; $i -lt $power.Count; $i++) {
    $r = $i + 2
    if ($r -eq 11) {
        $ws.Cells.Item($r, 3).Formula = "=6*7"
    } else {
        $ws.Cells.Item($r, 3).Value = $power[$i]
    }
    $ws.Cells.Item($r, 4).Value = $weekly[$i]
    $ws.Cells.Item($r, 5).Value = $timeuse[$i]
}

# Column B width
$ws.Columns.Item(2).ColumnWidth = 20.28515625

# Selection
$ws.Range("B11").Select()
